$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Keep a throwaway copy of the sheet around so the header row's original
#    (centered) formatting can be re-applied later without the engine
#    minting a brand new, unused cell style for it.
# ---------------------------------------------------------------------------
$ws.Copy($null, $ws)
$backup = $wb.Worksheets.Item($ws.Index + 1)

# ---------------------------------------------------------------------------
# 1. Rewrite columns E:F in place. They used to hold
#       E = OriginalStudiesNumbers (numeric: 500/600/700/800)
#       F = RecordsNumber          (numeric: 1000/1200/1400/1600)
#    They now hold a "locator / value" pair:
#       E = stdy_type_locators (text label of which metric)
#       F = stdy_type_values   (the metric's numeric value - the old E numbers)
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "stdy_type_locators"
$ws.Range("F1").Value = "stdy_type_values"

$ws.Range("E2").Value = "original_studies"
$ws.Range("F2").Value = 500

$ws.Range("E3").Value = "records_number"
$ws.Range("F3").Value = 600

$ws.Range("E4").Value = "fulltext_review"
$ws.Range("F4").Value = 700

$ws.Range("E5").Value = "total_record_number"
$ws.Range("F5").Value = 800

# ---------------------------------------------------------------------------
# 2. Normalize formatting: the data rows no longer carry the (visually
#    inert) alternate cell style, and the per-column format override is
#    dropped too - strip it from every column, then restore the header
#    row's centered alignment from the pristine backup copy.
# ---------------------------------------------------------------------------
$ws.Range("A100:XFD100").EntireColumn.ClearFormats()
$ws.Range("B3:C5").ClearContents()

$backup.Range("A1:I1").Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$backup.Delete()

# ---------------------------------------------------------------------------
# 3. The two metric columns that used to follow (fullTextReviewRecordsNumber,
#    totalRecordsNumber) are no longer needed - drop them entirely. This
#    shifts the trailing Prisma_Image column left into column G.
# ---------------------------------------------------------------------------
$ws.Range("G1:H1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 4. Reflect the edited region in the view (matches the author scrolling to
#    column C and leaving the E1:F5 block selected after editing it).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E1:F5").Select()
